$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header labels (shared between the "FV2404" and "FV2410" column groups).
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) -> "<name>_FV2404"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2404"
}

# Column K (11) stays "diff" - untouched.

# Columns L..U (12..21) -> "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2410"
}

# Turn the used range into an Excel Table ("Table1") now that headers carry
# their final names, so the table's column definitions pick them up.
$tableRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
